# Apply edits to the "Blue Line Block Info" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update sheet view: change selection to I10 (also clears the old scrolled
# topLeftCell position since the view re-centers around the new selection)
$ws.Activate()
$ws.Range("I10").Select()

# Rows 2..16: update D (length), E (grade), I (elevation, literal value now), J stays formula
# D column starts at 50 and decreases by 1 each row; E column equals (row-1) i.e. block number
for ($row = 2; $row -le 16; $row++) {
    $blockNum = $row - 1
    $ws.Cells.Item($row, 4).Value = 51 - $blockNum   # D column: 50,49,...,36
    $ws.Cells.Item($row, 5).Value = $blockNum        # E column: 1,2,...,15

    # I column: remove formula, set literal computed value 2.0 + blockNum/10
    $ws.Cells.Item($row, 9).Value = 2.0 + ($blockNum / 10.0)
}

# J column formulas remain as cumulative sums (already correct formulas in workbook);
# force recalculation so cached values match
$excel.Calculate()
